# Revert "Merge pull request #48 from LakeFishing/main"
# The merged change had fixed a typo in cell A3 (投籃 -> 投藍) and left the
# selection on G5; reverting restores the original text in A3 and moves the
# active selection back to G5 (matching the pre-revert cursor position saved
# in the file being reverted TO).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the single-character typo in A3 (投籃 -> 投藍).
$ws.Range("A3").Value = "投藍"

# Restore the saved selection/active cell to G5.
$ws.Range("G5").Select() | Out-Null
